$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bill Configuration")

# Update the "Gratuity Based on" option text (shared string swap: Guest Count -> Order Amount)
$ws.Range("B10").Value = "Order Amount"

# Update percentage / amount values
$ws.Range("B2").Value = 66.66
$ws.Range("B5").Value = 33.33
$ws.Range("B7").Value = 19.63
$ws.Range("B11").Value = 75.99
$ws.Range("B12").Value = 12
$ws.Range("B13").Value = 99
$ws.Range("B14").Value = 51

# Toggle ON/OFF switches
$ws.Range("B3").Value = "OFF"
$ws.Range("B6").Value = "OFF"
$ws.Range("B8").Value = "ON"
$ws.Range("B9").Value = "OFF"
$ws.Range("B16").Value = "OFF"
$ws.Range("B17").Value = "OFF"

# Update the active selection to match the saved view state
$ws.Range("B9").Select() | Out-Null
